# ServerDatabaseConfig sheet gains two new columns ("ConnectionTimeOut" and
# "QueryTimeOut") inserted between the existing "DatabaseName" and
# "WorkersCount" columns, with values 30 and 60 on the data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ServerDatabaseConfig")
$ws.Activate()

# Insert two blank columns at C:D, pushing the old WorkersCount column to E.
[void]$ws.Columns("C:D").Insert()

# New header row entries.
$ws.Range("C1").Value = "ConnectionTimeOut"
$ws.Range("D1").Value = "QueryTimeOut"

# New data row entries.
$ws.Range("C2").Value = 30
$ws.Range("D2").Value = 60

# Match the column widths of the surrounding "DatabaseName" column.
$ws.Columns("C:D").ColumnWidth = $ws.Columns("B:B").ColumnWidth

# Leave the same selection state recorded in the saved workbook.
[void]$ws.Range("D10").Select()
